$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4592192
$ws.Range("I33").Value = 6709665
$ws.Range("J33").Value = 4334
$ws.Range("K33").Value = 6709665
$ws.Range("L33").Value = 4334
$ws.Range("M33").Value = -6709436
$ws.Range("N33").Value = -4792

$ws.Range("H40").Value = 1739.1666
$ws.Range("I40").Value = 1604.8572
$ws.Range("J40").Value = 1927.2
$ws.Range("K40").Value = 1604.8572
$ws.Range("L40").Value = 1927.2
$ws.Range("M40").Value = -1429.8572
$ws.Range("N40").Value = -2277.2

$ws.Range("H64").Value = 9697.5
$ws.Range("I64").Value = 6930
$ws.Range("K64").Value = 6930
$ws.Range("M64").Value = -6682

$ws.Range("H67").Value = 9697.5
$ws.Range("I67").Value = 6930
$ws.Range("K67").Value = 6930
$ws.Range("M67").Value = -6072

$ws.Range("H70").Value = 1974.75
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1974.75
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5924.25
$ws.Range("N70").Value = -6464.25
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 1974.75
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1974.75
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5924.25
$ws.Range("N73").Value = -7796.25
$ws.Range("M73").ClearContents()

$ws.Range("H80").Value = 104167650
$ws.Range("I80").Value = 250000640
$ws.Range("K80").Value = 750001920
$ws.Range("M80").Value = -750000922

$ws.Range("H83").Value = 104167650
$ws.Range("I83").Value = 250000640
$ws.Range("K83").Value = 2250005760
$ws.Range("M83").Value = -2250000768

$ws.Range("H132").Value = 3429.4375
$ws.Range("I132").Value = 3619.3572
$ws.Range("K132").Value = 10858.0716
$ws.Range("M132").Value = -8328.071599999999

$ws.Range("H137").Value = 11112290
$ws.Range("I137").Value = 1202.1666
$ws.Range("K137").Value = 3606.4998
$ws.Range("M137").Value = -1056.4998

$ws.Range("H141").Value = 2019.4445
$ws.Range("I141").Value = 1710.7142
$ws.Range("K141").Value = 5132.142599999999
$ws.Range("M141").Value = 47.85740000000078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 713.3158
$ws.Range("I2").Value = 654.75
$ws.Range("K2").Value = 654.75
$ws.Range("M2").Value = -541.75

$ws.Range("H74").Value = 573639.7
$ws.Range("I74").Value = 2698.65
$ws.Range("J74").Value = 1174630.4
$ws.Range("K74").Value = 2698.65
$ws.Range("L74").Value = 1174630.4
$ws.Range("M74").Value = -1824.65
$ws.Range("N74").Value = -1176378.4

$ws.Range("H77").Value = 573639.7
$ws.Range("I77").Value = 2698.65
$ws.Range("J77").Value = 1174630.4
$ws.Range("K77").Value = 13493.25
$ws.Range("L77").Value = 5873152
$ws.Range("M77").Value = -9125.25
$ws.Range("N77").Value = -5881888

$ws.Range("H110").Value = 1129.6666
$ws.Range("I110").Value = 1143.6316
$ws.Range("K110").Value = 1143.6316
$ws.Range("M110").Value = 901.3684000000001

$ws.Range("H116").Value = 713.3158
$ws.Range("I116").Value = 654.75
$ws.Range("K116").Value = 654.75
$ws.Range("M116").Value = 1639.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 713.3158
$ws.Range("I3").Value = 654.75
$ws.Range("K3").Value = 654.75
$ws.Range("M3").Value = -540.75

$ws.Range("H7").Value = 997000
$ws.Range("I7").Value = 997000
$ws.Range("K7").Value = 997000
$ws.Range("M7").Value = -996887

$ws.Range("H50").Value = 80000
$ws.Range("J50").Value = 80000
$ws.Range("L50").Value = 80000
$ws.Range("N50").Value = -81148

$ws.Range("H105").Value = 9292.643
$ws.Range("I105").Value = 9116.5
$ws.Range("K105").Value = 9116.5
$ws.Range("M105").Value = -7369.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 15000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15776
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 15000
$ws.Range("N40").Value = -15320
$ws.Range("M40").ClearContents()

$ws.Range("H86").Value = 17724.25
$ws.Range("I86").Value = 5250
$ws.Range("J86").Value = 30198.5
$ws.Range("K86").Value = 5250
$ws.Range("L86").Value = 30198.5
$ws.Range("M86").Value = -4127
$ws.Range("N86").Value = -32444.5

$ws.Range("H89").Value = 17724.25
$ws.Range("I89").Value = 5250
$ws.Range("J89").Value = 30198.5
$ws.Range("K89").Value = 26250
$ws.Range("L89").Value = 150992.5
$ws.Range("M89").Value = -20634
$ws.Range("N89").Value = -162224.5

$ws.Range("H107").Value = 885.5599999999999
$ws.Range("I107").Value = 892.9
$ws.Range("J107").Value = 856.2
$ws.Range("K107").Value = 892.9
$ws.Range("L107").Value = 856.2
$ws.Range("M107").Value = 1027.1
$ws.Range("N107").Value = -4696.2

$ws.Range("H132").Value = 16719446
$ws.Range("I132").Value = 65065.5
$ws.Range("K132").Value = 195196.5
$ws.Range("M132").Value = -192666.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1121.5834
$ws.Range("I2").Value = 514.38464
$ws.Range("J2").Value = 1464.7826
$ws.Range("K2").Value = 3086.30784
$ws.Range("L2").Value = 8788.695599999999
$ws.Range("M2").Value = -2973.30784
$ws.Range("N2").Value = -9014.695599999999

$ws.Range("H9").Value = 75666664
$ws.Range("J9").Value = 75666664
$ws.Range("L9").Value = 226999992
$ws.Range("N9").Value = -227000440

$ws.Range("H34").Value = 2933.3125
$ws.Range("I34").Value = 618
$ws.Range("J34").Value = 3264.0715
$ws.Range("K34").Value = 1854
$ws.Range("L34").Value = 9792.2145
$ws.Range("M34").Value = -1770
$ws.Range("N34").Value = -9960.2145

$ws.Range("H38").Value = 521.5
$ws.Range("J38").Value = 826.2
$ws.Range("L38").Value = 2478.6
$ws.Range("N38").Value = -3172.6

$ws.Range("H113").Value = 267.16
$ws.Range("J113").Value = 258.05
$ws.Range("L113").Value = 774.1500000000001
$ws.Range("N113").Value = -5114.15

$ws.Range("H138").Value = 5624.933
$ws.Range("I138").Value = 5531.25
$ws.Range("J138").Value = 5999.6665
$ws.Range("K138").Value = 16593.75
$ws.Range("L138").Value = 17998.9995
$ws.Range("M138").Value = -11453.75
$ws.Range("N138").Value = -28278.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 127.64286
$ws.Range("I2").Value = 85.25
$ws.Range("J2").Value = 184.16667
$ws.Range("K2").Value = 85.25
$ws.Range("L2").Value = 184.16667
$ws.Range("M2").Value = 27.75
$ws.Range("N2").Value = -410.16667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2874.8333
$ws.Range("I9").Value = 3125
$ws.Range("K9").Value = 3125
$ws.Range("M9").Value = -2901

$ws.Range("H22").Value = 4660.0527
$ws.Range("J22").Value = 6397.3335
$ws.Range("L22").Value = 6397.3335
$ws.Range("N22").Value = -6987.3335

$ws.Range("H27").Value = 4660.0527
$ws.Range("J27").Value = 6397.3335
$ws.Range("L27").Value = 6397.3335
$ws.Range("N27").Value = -6611.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H96").Value = 2424.3333
$ws.Range("J96").Value = 1698
$ws.Range("L96").Value = 1698
$ws.Range("N96").Value = -4444

$ws.Range("H107").Value = 350.41177
$ws.Range("I107").Value = 383.65384
$ws.Range("K107").Value = 1150.96152
$ws.Range("M107").Value = 769.0384799999999

$ws.Range("H136").Value = 30457.176
$ws.Range("I136").Value = 35198.484
$ws.Range("J136").Value = 2957.6
$ws.Range("K136").Value = 105595.452
$ws.Range("L136").Value = 8872.799999999999
$ws.Range("M136").Value = -103045.452
$ws.Range("N136").Value = -13972.8

Write-Host "All changes applied"
